# PTD.xlsx - reportes horario y general PTD
# - Rename sheet Carga_Lectiva_PTD -> PTD
# - Update the remembered selection to the teacher-name box (B5:T5)
# - Top-align the comments/observations box (merged B93:AN95) so its
#   text starts at the top of the enlarged cells instead of the default
#   (bottom) vertical alignment.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet
$ws.Name = "PTD"

# 2. Leave the saved selection on B5:T5 (the docente name row)
$ws.Range("B5:T5").Select()

# 3. Vertical-align-top the merged comments box, keeping its existing
#    left horizontal alignment and borders intact
$ws.Range("B93:AN95").VerticalAlignment = -4160  # xlTop
